$wb = $excel.ActiveWorkbook

# --- 1. Build the two new "contingency" sheets from Line_Data ---------------
$lineData = $wb.Worksheets.Item("Line_Data")

# Line_Data_Con1: a copy of Line_Data placed right after it.
$lineData.Copy($null, $lineData)
$con1 = $wb.Worksheets.Item("Line_Data (2)")
$con1.Name = "Line_Data_Con1"

# Line_Data_Con2: a copy of Line_Data_Con1 placed right after it.
$con1.Copy($null, $con1)
$con2 = $wb.Worksheets.Item("Line_Data_Con1 (2)")
$con2.Name = "Line_Data_Con2"

# --- 2. Contingency 1: line 1-2 impedance doubled, susceptance halved -------
$con1.Range("C2").Formula = "=0.01938*2"
$con1.Range("D2").Formula = "=0.05917*2"
$con1.Range("E2").Formula = "=0.0528/2"

# --- 3. Contingency 2: line 1-5 removed from service ------------------------
$con2.Rows.Item(3).Delete()

# --- 4. View state: selections on each Line_Data* sheet ---------------------
$lineData.Range("A1:F18").Select()
$con1.Range("A1:F18").Select()
$con2.Range("F2").Select()

# Line_Data_Con2 ends up the active tab (mirrors the source workbook).
$con2.Activate()
